$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1 labels (plain, unstyled cells) ---
$ws.Range("C1").Value = "QCP"
$ws.Range("F1").Value = "Cardiac Arrest"
$ws.Range("H1").Value = "HumMod"

# --- Duplicate the formatted A2:D6 block into F2:I6 (new HumMod table) ---
$ws.Range("A2:D6").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 header labels (same captions as the first table)
$ws.Range("F2").Value = "Time"
$ws.Range("G2").Value = "Control"
$ws.Range("H2").Value = "30 Sec"
$ws.Range("I2").Value = "1 Min"

# Row labels (reuse existing shared strings) + HumMod "Control" column data
$ws.Range("F3").Value = "Blood Pressure(mmHg)"
$ws.Range("G3").Value = "120/79"

$ws.Range("F4").Value = "Cardiac Output(mL/min)"
$ws.Range("G4").Value = 5468

$ws.Range("F5").Value = "Ventilation(L/min)"
$ws.Range("G5").Value = 6.6

$ws.Range("F6").Value = "Symp Activity(Hz)"
$ws.Range("G6").Value = 1.5

# H3:I6 stay blank - the divide-by-zero error halted the run before 30 Sec / 1 Min data
$ws.Range("H3:I6").ClearContents()

# --- Explanatory note block, rows 8-10, merged A8:I10, yellow fill ---

# Step 1: stage the A8 "label" style at scratch cell K1 - same font/wrap/valign as the
# row-label cells (copied from A2), but with no border, yellow fill, and centered text.
$ws.Range("A2").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K1").Borders.LineStyle = -4142
$ws.Range("K1").Interior.Color = 65535
$ws.Range("K1").HorizontalAlignment = -4108

# Step 2: stage the plain block style at scratch cell K2 - default font, yellow fill, centered.
$ws.Range("K2").Interior.Color = 65535
$ws.Range("K2").HorizontalAlignment = -4108

# Step 3: apply the plain block style across the whole merged area first ...
$ws.Range("K2").Copy()
$ws.Range("A8:I10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 4: ... then overwrite A8 alone with the label style + the note text.
$ws.Range("K1").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A8").Value = "Patient does not reach 30 seconds before a divide by zero error in Structure\CO\CO.DES occurs The variable involved is Bronchi[CO] - (Uptake / Breathing.AlveolarVentilation(STPD))"

# clean up the scratch cells
$ws.Range("K1:K2").Clear()

$ws.Range("A8:I10").Merge()

$ws.Range("H3").Select()
